$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1220.1515
$ws.Cells.Item(17, 10).Value = 1220.1515
$ws.Cells.Item(17, 12).Value = 3660.4545
$ws.Cells.Item(17, 14).Value = -3996.4545
$ws.Cells.Item(28, 8).Value = 1221.6154
$ws.Cells.Item(28, 9).Value = 1353.1818
$ws.Cells.Item(28, 10).Value = 498
$ws.Cells.Item(28, 11).Value = 1353.1818
$ws.Cells.Item(28, 12).Value = 498
$ws.Cells.Item(28, 13).Value = -868.1818000000001
$ws.Cells.Item(28, 14).Value = -1468
$ws.Cells.Item(40, 8).Value = 1915.2307
$ws.Cells.Item(40, 9).Value = 1929.9
$ws.Cells.Item(40, 10).Value = 1866.3334
$ws.Cells.Item(40, 11).Value = 1929.9
$ws.Cells.Item(40, 12).Value = 1866.3334
$ws.Cells.Item(40, 13).Value = -1754.9
$ws.Cells.Item(40, 14).Value = -2216.3334
$ws.Cells.Item(48, 8).Value = 3019
$ws.Cells.Item(48, 10).Value = 3019
$ws.Cells.Item(48, 12).Value = 9057
$ws.Cells.Item(48, 14).Value = -9641
$ws.Cells.Item(56, 8).Value = 3019
$ws.Cells.Item(56, 10).Value = 3019
$ws.Cells.Item(56, 12).Value = 9057
$ws.Cells.Item(56, 14).Value = -10125
$ws.Cells.Item(98, 8).Value = 7416.0415
$ws.Cells.Item(98, 9).Value = 8378.904
$ws.Cells.Item(98, 10).Value = 676
$ws.Cells.Item(98, 11).Value = 8378.904
$ws.Cells.Item(98, 12).Value = 676
$ws.Cells.Item(98, 13).Value = -6880.904
$ws.Cells.Item(98, 14).Value = -3672
$ws.Cells.Item(118, 8).Value = 1004.875
$ws.Cells.Item(118, 9).Value = 434.14285
$ws.Cells.Item(118, 10).Value = 5000
$ws.Cells.Item(118, 11).Value = 1302.42855
$ws.Cells.Item(118, 12).Value = 15000
$ws.Cells.Item(118, 13).Value = 354.5714499999999
$ws.Cells.Item(118, 14).Value = -18314
$ws.Cells.Item(122, 8).Value = 7416.0415
$ws.Cells.Item(122, 9).Value = 8378.904
$ws.Cells.Item(122, 10).Value = 676
$ws.Cells.Item(122, 11).Value = 25136.712
$ws.Cells.Item(122, 12).Value = 2028
$ws.Cells.Item(122, 13).Value = -22686.712
$ws.Cells.Item(122, 14).Value = -6928
$ws.Cells.Item(129, 8).Value = 841.0417
$ws.Cells.Item(129, 9).Value = 291.16666
$ws.Cells.Item(129, 10).Value = 919.5952
$ws.Cells.Item(129, 11).Value = 873.4999799999999
$ws.Cells.Item(129, 12).Value = 2758.7856
$ws.Cells.Item(129, 13).Value = 4126.50002
$ws.Cells.Item(129, 14).Value = -12758.7856
$ws.Cells.Item(132, 8).Value = 9264048
$ws.Cells.Item(132, 9).Value = 12826154
$ws.Cells.Item(132, 10).Value = 2571.1
$ws.Cells.Item(132, 11).Value = 38478462
$ws.Cells.Item(132, 12).Value = 7713.299999999999
$ws.Cells.Item(132, 13).Value = -38475932
$ws.Cells.Item(132, 14).Value = -12773.3
$ws.Cells.Item(135, 8).Value = 964.75
$ws.Cells.Item(135, 9).Value = 286.15625
$ws.Cells.Item(135, 10).Value = 3679.125
$ws.Cells.Item(135, 11).Value = 2575.40625
$ws.Cells.Item(135, 12).Value = 33112.125
$ws.Cells.Item(135, 13).Value = -40.40625
$ws.Cells.Item(135, 14).Value = -38182.125
$ws.Cells.Item(138, 8).Value = 1409.77
$ws.Cells.Item(138, 9).Value = 842.23334
$ws.Cells.Item(138, 10).Value = 1653
$ws.Cells.Item(138, 11).Value = 2526.70002
$ws.Cells.Item(138, 12).Value = 4959
$ws.Cells.Item(138, 13).Value = 2613.29998
$ws.Cells.Item(138, 14).Value = -15239
$ws.Cells.Item(141, 8).Value = 571.05554
$ws.Cells.Item(141, 9).Value = 571.05554
$ws.Cells.Item(141, 11).Value = 1713.16662
$ws.Cells.Item(141, 13).Value = 3466.83338

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4123.492
$ws.Cells.Item(32, 9).Value = 3808.0894
$ws.Cells.Item(32, 11).Value = 3808.0894
$ws.Cells.Item(32, 13).Value = -3521.0894
$ws.Cells.Item(132, 8).Value = 1727.079
$ws.Cells.Item(132, 9).Value = 1682.3684
$ws.Cells.Item(132, 10).Value = 1771.7894
$ws.Cells.Item(132, 11).Value = 5047.1052
$ws.Cells.Item(132, 12).Value = 5315.3682
$ws.Cells.Item(132, 13).Value = -2517.1052
$ws.Cells.Item(132, 14).Value = -10375.3682

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 1115.8334
$ws.Cells.Item(22, 9).Value = 850
$ws.Cells.Item(22, 10).Value = 1248.75
$ws.Cells.Item(22, 11).Value = 850
$ws.Cells.Item(22, 12).Value = 1248.75
$ws.Cells.Item(22, 13).Value = -677
$ws.Cells.Item(22, 14).Value = -1594.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2498.8
$ws.Cells.Item(31, 9).Value = 2490
$ws.Cells.Item(31, 10).Value = 2504.6667
$ws.Cells.Item(31, 11).Value = 2490
$ws.Cells.Item(31, 12).Value = 2504.6667
$ws.Cells.Item(31, 13).Value = -2195
$ws.Cells.Item(31, 14).Value = -3094.6667
$ws.Cells.Item(34, 8).Value = 2498.8
$ws.Cells.Item(34, 9).Value = 2490
$ws.Cells.Item(34, 10).Value = 2504.6667
$ws.Cells.Item(34, 11).Value = 2490
$ws.Cells.Item(34, 12).Value = 2504.6667
$ws.Cells.Item(34, 13).Value = -2288
$ws.Cells.Item(34, 14).Value = -2908.6667
$ws.Cells.Item(107, 8).Value = 574.2308
$ws.Cells.Item(107, 9).Value = 550.8
$ws.Cells.Item(107, 10).Value = 606.1818
$ws.Cells.Item(107, 11).Value = 550.8
$ws.Cells.Item(107, 12).Value = 606.1818
$ws.Cells.Item(107, 13).Value = 1369.2
$ws.Cells.Item(107, 14).Value = -4446.1818
$ws.Cells.Item(132, 8).Value = 1519
$ws.Cells.Item(132, 9).Value = 859.5833
$ws.Cells.Item(132, 10).Value = 3497.25
$ws.Cells.Item(132, 11).Value = 2578.7499
$ws.Cells.Item(132, 12).Value = 10491.75
$ws.Cells.Item(132, 13).Value = -48.7498999999998
$ws.Cells.Item(132, 14).Value = -15551.75
$ws.Cells.Item(134, 8).Value = 956.8570999999999
$ws.Cells.Item(134, 9).Value = 650
$ws.Cells.Item(134, 10).Value = 1079.6
$ws.Cells.Item(134, 11).Value = 1950
$ws.Cells.Item(134, 12).Value = 3238.8
$ws.Cells.Item(134, 13).Value = 585
$ws.Cells.Item(134, 14).Value = -8308.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(132, 8).Value = 2015.6666
$ws.Cells.Item(132, 9).Value = 2200
$ws.Cells.Item(132, 11).Value = 19800
$ws.Cells.Item(132, 13).Value = -17270
$ws.Cells.Item(137, 8).Value = 22063866
$ws.Cells.Item(137, 9).Value = 44119040
$ws.Cells.Item(137, 10).Value = 8690.764999999999
$ws.Cells.Item(137, 11).Value = 132357120
$ws.Cells.Item(137, 12).Value = 26072.295
$ws.Cells.Item(137, 13).Value = -132352020
$ws.Cells.Item(137, 14).Value = -36272.295

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 2198.889
$ws.Cells.Item(126, 9).Value = 1881.6666
$ws.Cells.Item(126, 11).Value = 5644.9998
$ws.Cells.Item(126, 13).Value = -3174.9998
$ws.Cells.Item(132, 8).Value = 2462.5454
$ws.Cells.Item(132, 9).Value = 1582.3334
$ws.Cells.Item(132, 11).Value = 4747.0002
$ws.Cells.Item(132, 13).Value = -2217.0002

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2164.4285
$ws.Cells.Item(7, 9).Value = 2116.4167
$ws.Cells.Item(7, 10).Value = 2452.5
$ws.Cells.Item(7, 11).Value = 2116.4167
$ws.Cells.Item(7, 12).Value = 2452.5
$ws.Cells.Item(7, 13).Value = -2004.4167
$ws.Cells.Item(7, 14).Value = -2676.5
$ws.Cells.Item(22, 8).Value = 907.7778
$ws.Cells.Item(22, 9).Value = 692.5
$ws.Cells.Item(22, 10).Value = 969.2857
$ws.Cells.Item(22, 11).Value = 692.5
$ws.Cells.Item(22, 12).Value = 969.2857
$ws.Cells.Item(22, 13).Value = -397.5
$ws.Cells.Item(22, 14).Value = -1559.2857
$ws.Cells.Item(27, 8).Value = 907.7778
$ws.Cells.Item(27, 9).Value = 692.5
$ws.Cells.Item(27, 10).Value = 969.2857
$ws.Cells.Item(27, 11).Value = 692.5
$ws.Cells.Item(27, 12).Value = 969.2857
$ws.Cells.Item(27, 13).Value = -585.5
$ws.Cells.Item(27, 14).Value = -1183.2857
$ws.Cells.Item(46, 8).Value = 1533.3334
$ws.Cells.Item(46, 9).Value = 800
$ws.Cells.Item(46, 10).Value = 3000
$ws.Cells.Item(46, 11).Value = 800
$ws.Cells.Item(46, 12).Value = 3000
$ws.Cells.Item(46, 13).Value = -612
$ws.Cells.Item(46, 14).Value = -3376
$ws.Cells.Item(61, 8).Value = 1318.3334
$ws.Cells.Item(61, 9).Value = 300
$ws.Cells.Item(61, 10).Value = 1522
$ws.Cells.Item(61, 11).Value = 300
$ws.Cells.Item(61, 12).Value = 1522
$ws.Cells.Item(61, 13).Value = -98
$ws.Cells.Item(61, 14).Value = -1926
$ws.Cells.Item(113, 8).Value = 1318.3334
$ws.Cells.Item(113, 9).Value = 300
$ws.Cells.Item(113, 10).Value = 1522
$ws.Cells.Item(113, 11).Value = 300
$ws.Cells.Item(113, 12).Value = 1522
$ws.Cells.Item(113, 13).Value = 1870
$ws.Cells.Item(113, 14).Value = -5862
$ws.Cells.Item(126, 8).Value = 2164.4285
$ws.Cells.Item(126, 9).Value = 2116.4167
$ws.Cells.Item(126, 10).Value = 2452.5
$ws.Cells.Item(126, 11).Value = 6349.250100000001
$ws.Cells.Item(126, 12).Value = 7357.5
$ws.Cells.Item(126, 13).Value = -3879.250100000001
$ws.Cells.Item(126, 14).Value = -12297.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 47620320
$ws.Cells.Item(126, 9).Value = 52632852
$ws.Cells.Item(126, 10).Value = 1252.5
$ws.Cells.Item(126, 11).Value = 157898556
$ws.Cells.Item(126, 12).Value = 3757.5
$ws.Cells.Item(126, 13).Value = -157896086
$ws.Cells.Item(126, 14).Value = -8697.5
$ws.Cells.Item(132, 8).Value = 1583.5416
$ws.Cells.Item(132, 9).Value = 1300.35
$ws.Cells.Item(132, 10).Value = 2999.5
$ws.Cells.Item(132, 11).Value = 3901.05
$ws.Cells.Item(132, 12).Value = 8998.5
$ws.Cells.Item(132, 13).Value = -1371.05
$ws.Cells.Item(132, 14).Value = -14058.5
